# Regenerate save_data to use K (strikeouts) instead of Strike# for column G.
# Updates the "K" column values (G2:G34) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 2
    4  = 2
    5  = 8
    6  = 2
    7  = 10
    8  = 4
    9  = 1
    10 = 4
    11 = 5
    12 = 8
    13 = 4
    14 = 2
    15 = 11
    16 = 6
    17 = 3
    18 = 5
    19 = 8
    20 = 7
    21 = 6
    22 = 8
    23 = 7
    24 = 4
    25 = 3
    26 = 8
    27 = 4
    28 = 8
    29 = 6
    30 = 6
    31 = 5
    32 = 8
    33 = 6
    34 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
